# Update bitcoin_buys.xlsx after running on 2025-05-21
# Appends a new row (17) with that day's purchase data, matching the
# style of the most recently appended rows (plain text date, no special
# cell formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a literal date string like "05/18/2025" on recent rows
# (not an Excel date serial). Force text interpretation via NumberFormat
# so Excel doesn't auto-convert the string into a date serial number,
# then restore the default "Normal" style so no stray formatting is left
# on the cell (matching the unstyled cells used by the other recent rows).
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "05/21/2025"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = 0.0004661999999999999
$ws.Range("C17").Value = 107250.1072501073
$ws.Range("D17").Value = 50
